# Convert the year header labels (e.g. "1960 [YR1960]") in row 1 from text
# into plain numeric year values 1960..2019 across columns E:BL, and
# left-align them (they were left-aligned text, now they're left-aligned
# numbers) so the "line graph" code downstream can treat the header row as
# a numeric series for correlation/for-loop processing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 5   # column E
$lastCol  = 64  # column BL
$startYear = 1960

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $year = $startYear + ($col - $firstCol)
    $ws.Cells.Item(1, $col).Value = $year
}

# Match the original left-aligned look of the text year labels.
$ws.Range("E1:BL1").HorizontalAlignment = -4131   # xlHAlignLeft

# Leave the header-year range selected, as in the authored workbook.
[void]$ws.Range("E1:BL1").Select()
